# Adding TestCases for IAM
# Appends four new rows (TestCase_B64 .. TestCase_B67) to the "Test Cases"
# sheet, continuing directly after the existing TestCase_B63 row (row 64).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy the formatting of the last existing data row (row 64) down onto the
# four new rows, so the new cells pick up the same borders/fills/wrap
# settings already used throughout this block of the sheet.
$ws.Range("A64:E64").Copy()
$ws.Range("A65:E68").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    @("TestCase_B64", "OPQA-557", "Verify that the searched keyword present in the search text box doesn't change if any other content type is selected in the search drop down", "Y", "SKIP"),
    @("TestCase_B65", "OPQA-386", "Verify that the searched keyword doesn't change in the search text box if any other content type is selected in the left navigation pane", "Y", "SKIP"),
    @("TestCase_B66", "OPQA-387", "Verify that counts of search results of all the content types should get displayed irrespective of the content type chosen for searching", "Y", "SKIP"),
    @("TestCase_B67", "OPQA-263", "Verify that ALL search results count is equal to the count of search results of other content types(ARTICLES+PATENTS+POSTS+PEOPLE)", "Y", "PASS")
)

$startRow = 65
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}

$ws.Range("D10").Select()
